{"js": "const table = context.document.body.tables.getFirst();\n\nconst updates = [\n  { row: 0, col: 0, text: \"27\u00f76=4, 3\" },\n  { row: 0, col: 1, text: \"21\u00f76=3, 3\" },\n  { row: 0, col: 2, text: \"89\u00f79=9, 8\" },\n  { row: 0, col: 3, text: \"47\u00f76=7, 5\" },\n  { row: 0, col: 4, text: \"13\u00f77=1, 6\" },\n  { row: 4, col: 0, text: \"24\u00f76=4, 0\" },\n  { row: 4, col: 1, text: \"91\u00f76=15, 1\" },\n  { row: 4, col: 2, text: \"55\u00f78=6, 7\" },\n  { row: 4, col: 3, text: \"99\u00f76=16, 3\" },\n  { row: 4, col: 4, text: \"31\u00f78=3, 7\" },\n  { row: 8, col: 0, text: \"70\u00f74=17, 2\" },\n  { row: 8, col: 1, text: \"53\u00f77=7, 4\" },\n  { row: 8, col: 2, text: \"83\u00f76=13, 5\" },\n  { row: 8, col: 3, text: \"88\u00f72=44, 0\" },\n  { row: 8, col: 4, text: \"96\u00f72=48, 0\" },\n  { row: 12, col: 0, text: \"12\u00f77=1, 5\" },\n  { row: 12, col: 1, text: \"17\u00f79=1, 8\" },\n  { row: 12, col: 2, text: \"87\u00f78=10, 7\" },\n  { row: 12, col: 3, text: \"33\u00f76=5, 3\" },\n  { row: 12, col: 4, text: \"45\u00f72=22, 1\" },\n  { row: 16, col: 0, text: \"33\u00f79=3, 6\" },\n  { row: 16, col: 1, text: \"69\u00f76=11, 3\" },\n  { row: 16, col: 2, text: \"47\u00f78=5, 7\" },\n  { row: 16, col: 3, text: \"75\u00f72=37, 1\" },\n  { row: 16, col: 4, text: \"64\u00f76=10, 4\" },\n];\n\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  const rng = cell.body.getRange();\n  rng.insertText(u.text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$values = @(\n    @(1, 1, \"27\u00f76=4, 3\"),\n    @(1, 2, \"21\u00f76=3, 3\"),\n    @(1, 3, \"89\u00f79=9, 8\"),\n    @(1, 4, \"47\u00f76=7, 5\"),\n    @(1, 5, \"13\u00f77=1, 6\"),\n    @(5, 1, \"24\u00f76=4, 0\"),\n    @(5, 2, \"91\u00f76=15, 1\"),\n    @(5, 3, \"55\u00f78=6, 7\"),\n    @(5, 4, \"99\u00f76=16, 3\"),\n    @(5, 5, \"31\u00f78=3, 7\"),\n    @(9, 1, \"70\u00f74=17, 2\"),\n    @(9, 2, \"53\u00f77=7, 4\"),\n    @(9, 3, \"83\u00f76=13, 5\"),\n    @(9, 4, \"88\u00f72=44, 0\"),\n    @(9, 5, \"96\u00f72=48, 0\"),\n    @(13, 1, \"12\u00f77=1, 5\"),\n    @(13, 2, \"17\u00f79=1, 8\"),\n    @(13, 3, \"87\u00f78=10, 7\"),\n    @(13, 4, \"33\u00f76=5, 3\"),\n    @(13, 5, \"45\u00f72=22, 1\"),\n    @(17, 1, \"33\u00f79=3, 6\"),\n    @(17, 2, \"69\u00f76=11, 3\"),\n    @(17, 3, \"47\u00f78=5, 7\"),\n    @(17, 4, \"75\u00f72=37, 1\"),\n    @(17, 5, \"64\u00f76=10, 4\"),\n)\n\nforeach ($v in $values) {\n    $cell = $t.Cell($v[0], $v[1])\n    $cell.Range.Text = $v[2]\n}\n"}
